$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27: A27 changes from "Q013" to the new code "Q014" ---
$ws.Range("A27").Value = "Q014"

# --- Row 28 (new) ---
$ws.Range("A28").Value = "A001"
$ws.Range("B28").Value = "F"
# C28 must be stored as a genuine numeric literal (like C1) even though the
# column is formatted as Text (numFmtId 49, "@"). Round-trip the number
# format so the COM layer doesn't coerce the value to a string, then put
# the cell's format back to its original Text format.
$ws.Range("C28").NumberFormat = "General"
$ws.Range("C28").Value = 21
$ws.Range("C28").NumberFormat = "@"
$ws.Range("D28").Value = "001"
$ws.Range("E28").Value = "Underweight"
$ws.Range("F28").Value = "12"
$ws.Range("G28").Value = "1-1-1996"

# --- Row 29 (new) ---
$ws.Range("A29").Value = "A $$$ 283"
$ws.Range("B29").Value = "Ma le"
$ws.Range("C29").Value = "2$0"
$ws.Range("D29").Value = "   444"
$ws.Range("E29").Value = "Obese"
$ws.Range("F29").Value = "12"
$ws.Range("G29").Value = "1st Feb 1998"

# --- Row 30 (new) ---
$ws.Range("A30").Value = "A222"
$ws.Range("B30").Value = "Girl"
$ws.Range("C30").Value = "20"
$ws.Range("D30").Value = "777"
$ws.Range("E30").Value = "Normal"
$ws.Range("F30").Value = "12"
$ws.Range("G30").Value = "2nd January 1998"

# Put the cursor on the last entered cell, matching the saved selection.
$ws.Range("A30").Select() | Out-Null
